# Fills in the three "half-year summary" comment cells in the certificate
# table (נביא / הנדסה / חשבון rows) that were previously left blank.
#
# Each subject row is a 1x2 table: the first cell holds the subject label
# (e.g. "נביא"), the second cell is the free-text comment that needs to be
# populated. We locate each row by its label text so the script is not
# dependent on a brittle, hard-coded table index.

$d = $word.ActiveDocument

function Set-CommentForSubject($label, $text) {
    for ($i = 1; $i -le $d.Tables.Count; $i++) {
        $tbl = $d.Tables.Item($i)
        $labelText = $tbl.Cell(1, 1).Range.Text
        if ($labelText -match $label) {
            $r = $tbl.Cell(1, 2).Range
            # Trim the trailing cell-mark/paragraph-mark characters so we
            # only replace the visible text, keeping the existing run /
            # paragraph formatting intact.
            $r.End = $r.End - 1
            $r.Text = $text
            return
        }
    }
}

Set-CommentForSubject "נביא" "במחצית זאת למדנו את ספר שמואל, עם מפרשים והרחבנו בהרבה מדרשים,`nהייתה אוירה קדושה,והנאה בלימוד.`nחיה את מיוחדת, בהצלחה!"

Set-CommentForSubject "הנדסה" "למדנו במחצית זאת על סוגי המצולעים השונים. כגון ריבוע, משולש, משושה, מעוין וכו'. למדנו כיצד מזהים כל מצולע ומה תכונותיו.`nחיה, את תלמידה מדהימה! המשיכי להצליח בדרכך!"

Set-CommentForSubject "חשבון" "במחצית זאת חזרנו את פעולות חשבון, חיבור, חיסור, כפל וחילוק,`nהתמקדנו בעיקר על לוח הכפל, פיתחנו שיטות לימוד רציניות,ולמדנו איך לזכור דברים בע`"פ.`nחיה את ילדה מקסימה, יש לך ראש חזק, שיהיה בהצלחה!"
